$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119; this shifts the previous rows
# 119-221 down to 120-222 and updates the sheet dimension accordingly.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44669
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = "Fruta"
$ws.Cells.Item(119, 7).Value = 100103
$ws.Cells.Item(119, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(119, 9).Value = 100103002
$ws.Cells.Item(119, 10).Value = "Ciruela"
$ws.Cells.Item(119, 11).Value = "Angeleno"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 130
$ws.Cells.Item(119, 14).Value = 10000
$ws.Cells.Item(119, 15).Value = 12000
$ws.Cells.Item(119, 16).Value = 11231
$ws.Cells.Item(119, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(119, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(119, 19).Value = 624
$ws.Cells.Item(119, 20).Value = 18
